$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 47-48, pushing the existing rows 47:68 down to 49:70.
$ws.Range("A47:R48").Insert()

# New row 47: "Especial" quality entry dated 2023-05-03 (serial 45049)
$ws.Cells.Item(47,1).Value = 11
$ws.Cells.Item(47,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(47,3).Value = "Bíobío"
$ws.Cells.Item(47,4).Value = 45049
$ws.Cells.Item(47,5).Value = 8
$ws.Cells.Item(47,6).Value = 100112043
$ws.Cells.Item(47,7).Value = "Pepino dulce"
$ws.Cells.Item(47,8).Value = "Cultivar IV Región"
$ws.Cells.Item(47,9).Value = "Especial"
$ws.Cells.Item(47,10).Value = 50
$ws.Cells.Item(47,11).Value = 15000
$ws.Cells.Item(47,12).Value = 15000
$ws.Cells.Item(47,13).Value = 15000
$ws.Cells.Item(47,14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(47,15).Value = "Provincia de Limarí"
$ws.Cells.Item(47,16).Value = 833
$ws.Cells.Item(47,17).Value = 18
$ws.Cells.Item(47,18).Value = "Hortaliza"

# New row 48: "Primera" quality entry dated 2023-05-03 (serial 45049)
$ws.Cells.Item(48,1).Value = 11
$ws.Cells.Item(48,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(48,3).Value = "Bíobío"
$ws.Cells.Item(48,4).Value = 45049
$ws.Cells.Item(48,5).Value = 8
$ws.Cells.Item(48,6).Value = 100112043
$ws.Cells.Item(48,7).Value = "Pepino dulce"
$ws.Cells.Item(48,8).Value = "Cultivar IV Región"
$ws.Cells.Item(48,9).Value = "Primera"
$ws.Cells.Item(48,10).Value = 50
$ws.Cells.Item(48,11).Value = 13000
$ws.Cells.Item(48,12).Value = 13000
$ws.Cells.Item(48,13).Value = 13000
$ws.Cells.Item(48,14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(48,15).Value = "Provincia de Limarí"
$ws.Cells.Item(48,16).Value = 722
$ws.Cells.Item(48,17).Value = 18
$ws.Cells.Item(48,18).Value = "Hortaliza"
